$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Philip" / ID header cells in row 1.
$ws.Range("C1:E1").ClearContents()

# --- Fall 2022 / Spring 2022 / Summer 2022 section (rows 4-9) ---
# Row 4: POLS 1101 stays; Spring course changes from CPSC 3165 to CPSC 4148;
# a new Summer course (CPSC 4176) is added.
$ws.Range("C4").Value = "CPSC 4148"
$ws.Range("E4").Value = "CPSC 4176"
$ws.Range("F4").Value = 3

# Row 5: Spring course changes from CPSC 4135 to CPSC 4155.
$ws.Range("C5").Value = "CPSC 4155"

# Row 6: Spring course changes from CPSC 4148 to CPSC 4157.
$ws.Range("C6").Value = "CPSC 4157"

# Row 7: now holds CPSC 3165 (Fall) / CPSC 4175 (Spring) instead of the old
# CPSC 4000 total-credit row.
$ws.Range("A7").Value = "CPSC 3165"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4175"
$ws.Range("D7").Value = 3

# Row 8 (new): CPSC 4135.
$ws.Range("A8").Value = "CPSC 4135"
$ws.Range("B8").Value = 3

# Row 9 (new): CPSC 4000, 0 credits (was previously row 7).
$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0

# --- Remove the Fall/Spring/Summer 2023 data rows (13-15) ---
$ws.Range("A13:F15").ClearContents()

# --- Remove the entire Fall/Spring/Summer 2024 section (rows 21 and 29) ---
$ws.Range("A21:F21").ClearContents()
$ws.Range("A29:F29").ClearContents()
